$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Password (B2) and ConfirmPassword (D2) values
$ws.Range("B2").Value = "pavan@1236"
$ws.Range("D2").Value = "pavan@1236"

# Add hyperlinks to the new password values (mirrors existing email hyperlink pattern)
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:pavan@1236")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:pavan@1236")

# Re-apply the Hyperlink style so it matches the existing style (reuse, not duplicate)
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("D2").Style = "Hyperlink"

# Update the active selection
$ws.Range("B11").Select()
